$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new columns at F:G (everything from F onward shifts right by 2)
$ws.Columns("F:G").Insert()

# 2. Clear existing (now stale-positioned) hyperlinks; we'll re-add them at the right spots
$ws.Hyperlinks.Delete()

# 3. Header row (row 1)
$ws.Range("F1").Value = "Login_UserNameOnline"
$ws.Range("G1").Value = "Login_PasswordOnline"
$ws.Range("K1").Value = "PIDIntegation"

# 4. Row 2 - bring D2/E2 back to plain creds, and populate the new F2/G2 pair
#    First copy E2's (non-hyperlink) formatting onto D2 so D2 loses its old hyperlink style.
$ws.Range("E2").Copy($ws.Range("D2"))
$ws.Range("D2").Value = "ssharma"
$ws.Range("E2").Value = "Rockstar1"

#    Give F2 the "hyperlink" look by copying it from the cell that used to carry the
#    hyperlink style (D2 still has style 16 cached on the clipboard source... use a fresh copy)
$ws.Range("D2").Copy($ws.Range("F2"))
$ws.Range("F2").Value = "randycoplin2011@gmail.com"
$ws.Range("G2").Value = "password"
